$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("iconResource" and following shift right by one)
$ws.Columns("D").Insert()

# New header for the inserted column
$ws.Range("D1").Value = "invokeType"

# New trailing headers
$ws.Range("H1").Value = "effect.effectType"
$ws.Range("I1").Value = "effect.value"

# invokeType values mirror the "type" column (C)
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 2
# Normalize the new column's formatting uniformly (rows 4/7 would
# otherwise inherit the wrap-text style used by the "magic_sword" rows)
$ws.Range("D2:D7").WrapText = $false

# effect.effectType values mirror the "type" column as well
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 2
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 2

# effect.value is a simple running index
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 2
$ws.Range("I4").Value = 3
$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 6

$ws.Range("J6").Select() | Out-Null
